$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B11 used to show the (now unused) "R40" label; it is retyped as the
# literal text "1". Assigning the plain string "1" directly would be
# auto-coerced to the number 1 by Value's usual type inference, so we
# route the text through a formula first (which always yields a true
# text result for a quoted literal) and then flatten it back down to a
# plain value via copy / paste-special-values. That keeps the cell's
# existing style untouched while landing a genuine text cell.
$ws.Range("B11").Formula = '="1"'
$ws.Range("B11").Copy()
$ws.Range("B11").PasteSpecial(-4163)
$excel.CutCopyMode = $false
